# UsersCreationDetails.xlsx - "Completed User Creation part Exam Center"
#
# STAGE sheet: fill in Location (G), ID (D) and password/Tenant-like (E)
# details for the fpkcontroller row, and blank-out the other two rows'
# previously-placeholder ID/Tenant values to "Null". Also mark STAGE as
# the active/selected tab (it previously was LMSProd).
#
# LMSProd sheet: fill in the Location (G) column for the first two data
# rows, and stop being the selected tab.

$wb = $excel.ActiveWorkbook

$wsStage = $wb.Worksheets.Item("STAGE")
$wsLms   = $wb.Worksheets.Item("LMSProd")

# --- STAGE ---------------------------------------------------------------

# G2/G3 pick up the same "filled-in location" look already used by G4/G5
# (Arial font style) instead of the plain default border-only style.
$wsStage.Range("G4").Copy()
$wsStage.Range("G2").PasteSpecial(-4122)
$wsStage.Range("G3").PasteSpecial(-4122)
$wsStage.Range("G2").Value = "Canada 05"
$wsStage.Range("G3").Value = "Canada 05"

$wsStage.Range("D3").Value = "fpkcontroller"

# "159533" looks like a number - force it to stay text (matches the
# existing General-formatted style already on the cell) via a text
# formula that gets collapsed back down to a plain literal by
# Paste-Special Values (avoids creating a new number-format style).
$wsStage.Range("E3").Formula = '="159533"'
$wsStage.Range("E3").Copy()
$wsStage.Range("E3").PasteSpecial(-4163)

$wsStage.Range("D4").Value = "Null"
$wsStage.Range("E4").Value = "Null"

$wsStage.Range("D5").Value = "Null"
$wsStage.Range("E5").Value = "Null"

# --- LMSProd ---------------------------------------------------------------

# G2 picks up the filled-in-location style too (G3 already has it).
$wsLms.Range("G4").Copy()
$wsLms.Range("G2").PasteSpecial(-4122)
$wsLms.Range("G2").Value = "Japan1232"
$wsLms.Range("G3").Value = "Japan1232"

# --- selection / active tab ------------------------------------------------
# LMSProd loses tabSelected + its selection moves from G5 to G8.
$wsLms.Range("G8").Select()

# STAGE becomes the selected tab (was LMSProd) with G2 selected.
$wsStage.Range("G2").Select()
